$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Lesson #33 (row 37) and #34 (row 38): YouTube recordings uploaded ---
# Add the new hyperlinks first (this also writes the shared-string text),
# then restore the original "no special hyperlink theme" cell formatting by
# copying the format from a sibling cell that already carries a plain
# (non-hyperlink-themed) look, matching how the rest of the sheet's
# YouTube-link cells are styled.

$ws.Hyperlinks.Add($ws.Range("F37"), "https://youtu.be/dbPUQE2NsLo", [Type]::Missing, [Type]::Missing, "https://youtu.be/dbPUQE2NsLo")
$ws.Range("F37").Value = "https://youtu.be/dbPUQE2NsLo "

$ws.Hyperlinks.Add($ws.Range("F38"), "https://youtu.be/IAOZ1-2VPuQ", [Type]::Missing, [Type]::Missing, "https://youtu.be/IAOZ1-2VPuQ")
$ws.Range("F38").Value = "https://youtu.be/IAOZ1-2VPuQ "

# Restore the plain (non-auto-hyperlink-themed) cell format used throughout
# the rest of the "YouTube link" column.
$ws.Range("F29").Copy()
$ws.Range("F37").PasteSpecial(-4122)
$ws.Range("F29").Copy()
$ws.Range("F38").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row heights tighten up slightly now that the cells hold content.
$ws.Rows(37).RowHeight = 14.3
$ws.Rows(38).RowHeight = 14.3

# --- Clear the stray "amount (hours)" = 2 values that were placeholders
#     for not-yet-scheduled lessons (rows 39-44), matching how row 37/38
#     now carry real data while the rest of the table remains blank. ---
$ws.Range("D39").ClearContents()
$ws.Range("D40").ClearContents()
$ws.Range("D41").ClearContents()
$ws.Range("D42").ClearContents()
$ws.Range("D43").ClearContents()
$ws.Range("D44").ClearContents()

# Row 44 gains the same (empty) "Lesson name" style as its neighbours above.
$ws.Range("C43").Copy()
$ws.Range("C44").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 44's "amount" cell now shares the same format as the other blank rows.
$ws.Range("D43").Copy()
$ws.Range("D44").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Selection moves from E39 to D39 ---
$ws.Range("D39").Select() | Out-Null
